$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.392.14'
$ws.Range("E2").Value = '  +0.99%  '
$ws.Range("D3").Value = '3.161.62'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.580'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.44%  '
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("E10").Value = '  -0.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '3.716.94'
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = '64.451.31'
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").Value = '3.162.69'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '406.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.66'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.485'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("E25").Value = '  -2.78%  '
$ws.Range("E26").Value = '  -2.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("E30").Value = '  -1.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.87%  '
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.69'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("D37").Value = '2.680.19'
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.696'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '290.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0257'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.49'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.879'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.73%  '
